$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.255.77"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "2.060.79"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.93"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +2.73%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.69"
$ws.Range("E8").Value = "  +4.55%  "
$ws.Range("E9").Value = "  +2.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.06"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0758"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "2.363.88"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.34"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.76"
$ws.Range("E15").Value = "  +2.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.773"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.14"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "2.056.50"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "37.210.62"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.37"
$ws.Range("E20").Value = "  +14.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.07"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").Value = "0.0₃0810"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "224.91"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.36"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("E28").Value = "  +7.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.78"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.05"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.46"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0613"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.53"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.52"
$ws.Range("E36").Value = "  +6.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.74"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.26"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.56"
$ws.Range("E41").Value = "  +12.47%  "
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").Value = "1.489.94"
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.48"
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("E45").Value = "  +4.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0923"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.32"
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.14"
$ws.Range("E50").Value = "  +3.98%  "
$ws.Range("E51").Value = "  +2.35%  "
